# Insert a new weekly record as row 7 on the active sheet, pushing the
# existing rows 7..71 down to 8..72 (dimension grows from A1:R71 to A1:R72).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("7:7").Insert()

$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 45022
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112032
$ws.Cells.Item(7, 7).Value = "Zapallo italiano"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 5000
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = 5500
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 92
$ws.Cells.Item(7, 17).Value = 60
$ws.Cells.Item(7, 18).Value = "Hortaliza"
